$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10, pushing existing rows 10-20 down to 11-21.
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with this week's record (same as the old row 10
# except for the date and volume, which are this week's new figures).
$ws.Cells.Item(10, 1).Value = 2
$ws.Cells.Item(10, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44413
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 100112022
$ws.Cells.Item(10, 7).Value = "Arveja Verde"
$ws.Cells.Item(10, 8).Value = "Perfection"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 700
$ws.Cells.Item(10, 11).Value = 26000
$ws.Cells.Item(10, 12).Value = 28000
$ws.Cells.Item(10, 13).Value = 27000
$ws.Cells.Item(10, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 1080
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
